$d = $word.ActiveDocument

$replacements = @(
    @("30÷6=5, 0", "13÷9=1, 4"),
    @("50÷7=7, 1", "47÷4=11, 3"),
    @("40÷2=20, 0", "61÷7=8, 5"),
    @("17÷4=4, 1", "41÷7=5, 6"),
    @("15÷7=2, 1", "53÷4=13, 1"),
    @("31÷4=7, 3", "76÷3=25, 1"),
    @("85÷2=42, 1", "80÷7=11, 3"),
    @("50÷2=25, 0", "81÷2=40, 1"),
    @("20÷3=6, 2", "24÷8=3, 0"),
    @("87÷7=12, 3", "58÷2=29, 0"),
    @("70÷4=17, 2", "80÷9=8, 8"),
    @("98÷6=16, 2", "24÷4=6, 0"),
    @("63÷4=15, 3", "50÷4=12, 2"),
    @("65÷5=13, 0", "55÷2=27, 1"),
    @("94÷2=47, 0", "89÷4=22, 1"),
    @("10÷4=2, 2", "93÷5=18, 3"),
    @("31÷3=10, 1", "55÷7=7, 6"),
    @("37÷8=4, 5", "48÷5=9, 3"),
    @("96÷6=16, 0", "73÷9=8, 1"),
    @("32÷7=4, 4", "11÷4=2, 3"),
    @("47÷8=5, 7", "85÷4=21, 1"),
    @("12÷8=1, 4", "26÷3=8, 2"),
    @("83÷7=11, 6", "44÷2=22, 0"),
    @("90÷2=45, 0", "97÷8=12, 1"),
    @("18÷6=3, 0", "28÷7=4, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
